$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.302.55'
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = '1.844.05'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9982'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.56'
$ws.Range("E5").Value = '  +0.09%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6275'
$ws.Range("E6").Value = '  -0.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9999'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07458'
$ws.Range("E8").Value = '  -2.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2896'
$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.31'
$ws.Range("E10").Value = '  -2.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07719'
$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").Value = '1.843.99'
$ws.Range("E12").Value = '  -2.34%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.991'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6783'
$ws.Range("E14").Value = '  -0.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001011'
$ws.Range("E15").Value = '  -4.64%  '

$ws.Range("E16").Value = '  -1.78%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.148'
$ws.Range("E17").Value = '  -0.84%  '

$ws.Range("D18").Value = '29.294.12'
$ws.Range("E18").Value = '  -0.95%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.55'
$ws.Range("E19").Value = '  -0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.28'
$ws.Range("E20").Value = '  -0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.402'
$ws.Range("E22").Value = '  -0.78%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").Value = '  -0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.70'
$ws.Range("E24").Value = '  +0.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1375'
$ws.Range("E25").Value = '  -0.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.405'
$ws.Range("E26").Value = '  -0.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.57'
$ws.Range("E27").Value = '  -0.97%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06341'
$ws.Range("E28").Value = '  +12.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.394'
$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.473'
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.083'
$ws.Range("E31").Value = '  -1.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.051'
$ws.Range("E32").Value = '  -0.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.819'
$ws.Range("E33").Value = '  -1.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.143'
$ws.Range("E34").Value = '  -2.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6969'
$ws.Range("E35").Value = '  -0.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.581'
$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("E37").Value = '  +3.51%  '

$ws.Range("D38").Value = '1.243.36'
$ws.Range("E38").Value = '  +0.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01815'
$ws.Range("E39").Value = '  +0.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.531'
$ws.Range("E40").Value = '  +0.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9091'
$ws.Range("E41").Value = '  +0.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9981'
$ws.Range("E42").Value = '  -0.28%  '

$ws.Range("D43").Value = '2.002.69'
$ws.Range("E43").Value = '  -14.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.34'
$ws.Range("E44").Value = '  -1.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.24'
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.058'
$ws.Range("E46").Value = '  -2.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1165'
$ws.Range("E47").Value = '  +1.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.049'
$ws.Range("E48").Value = '  +0.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3948'
$ws.Range("E49").Value = '  -2.11%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.664'
$ws.Range("E50").Value = '  -0.96%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000114'
$ws.Range("E51").Value = '  -3.37%  '
